$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 511, shifting existing rows 511:624 down to 512:625
$ws.Rows.Item(511).Insert()

# Populate the newly inserted row 511 with the new data record
$ws.Cells.Item(511, 1).Value = 9
$ws.Cells.Item(511, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(511, 3).Value = "Metropolitana"
$ws.Cells.Item(511, 4).Value = 44711
$ws.Cells.Item(511, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(511, 5).Value = 13
$ws.Cells.Item(511, 6).Value = 100112040
$ws.Cells.Item(511, 7).Value = "Cilantro"
$ws.Cells.Item(511, 8).Value = "Sin especificar"
$ws.Cells.Item(511, 9).Value = "Primera"
$ws.Cells.Item(511, 10).Value = 70
$ws.Cells.Item(511, 11).Value = 9000
$ws.Cells.Item(511, 12).Value = 10000
$ws.Cells.Item(511, 13).Value = 9500
$ws.Cells.Item(511, 14).Value = "`$/docena de atados"
$ws.Cells.Item(511, 15).Value = "Región Metropolitana"
$ws.Cells.Item(511, 16).Value = 3167
$ws.Cells.Item(511, 17).Value = 3
$ws.Cells.Item(511, 18).Value = "Hortaliza"
